$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Shift existing columns D:E (old "Email" link columns) right by one
#    so a new column D ("FB Msg") can be inserted: old E -> F, old D -> E.
#    Destination ranges are cleared first (Copy onto a non-blank cell
#    from a blank source does not blank it back out in this engine),
#    then values + formats are copied across.
# ------------------------------------------------------------------
$ws.Range("F1:G18").Clear()
$ws.Range("E1:F18").Copy($ws.Range("F1"))

$ws.Range("E1:E18").Clear()
$ws.Range("D1:D18").Copy($ws.Range("E1"))

# 2. Wipe the old column D completely - we rebuild it from scratch.
$ws.Range("D1:D18").Clear()

Write-Host "shift done"
